# Update the "Title" paragraph style (and its linked "Title Char"
# character style) so the title font changes from Arial to Cambria,
# leaving the complex-script font (cs="Arial") untouched.

$d = $word.ActiveDocument

$titleStyle = $d.Styles("Title")
$titleStyle.Font.Name = "Cambria"

$titleCharStyle = $d.Styles("TitleChar")
$titleCharStyle.Font.Name = "Cambria"

Write-Output "Title style font -> $($titleStyle.Font.Name)"
Write-Output "TitleChar style font -> $($titleCharStyle.Font.Name)"
